$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 3 for the columns that differ
# (A, B, E, F, G, H, Q, R) plus the quirky empty "AF" marker cell which
# moves from row 3 to row 2.

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr2 = "$col`2"
    $addr3 = "$col`3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value = $v3
    $ws.Range($addr3).Value = $v2
}

# Move the empty-text marker cell from AF3 to AF2.
# A direct `.Value = ""` assignment clears/removes a cell entirely (same as
# real Excel) instead of leaving a present-but-empty text cell behind, so
# force an empty *text* value via the leading-apostrophe quote-prefix, then
# reset the style back to Normal so no stray quote-prefix formatting lingers.
$ws.Range("AF2").Value = "'"
$ws.Range("AF2").Style = "Normal"
$ws.Range("AF3").Value = $null
